$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Actividad para identificar los versos en un poema." is split
# into two runs ("...en" / " un poema") and loses its trailing period.
# ---------------------------------------------------------------------------
$found1 = $d.Content.Find.Execute("Actividad para identificar los versos en un poema.")
if ($found1) {
    $full1 = $d.Content
    $full1.Find.Execute("Actividad para identificar los versos en un poema.") | Out-Null
    $s1 = $full1.Start
    $e1 = $full1.End
    $splitAt = $s1 + "Actividad para identificar los versos en".Length
    $tail = $d.Range($splitAt, $e1)
    # Toggle formatting so the engine keeps the new text as a distinct run
    # instead of silently re-merging it with its identically formatted
    # neighbour.
    $tail.Font.Bold = 1
    $tail.Text = " un poema"
    $tail.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# Change 2: the keyword line becomes a quoted, spell-checked literal and the
# extra spaces after the commas are dropped; the run layout is rebuilt from
# scratch (quote / spellStart / Poema / ,versos, / palabras / , / literatura
# / spellEnd / quote) and the _GoBack bookmark moves here.
# ---------------------------------------------------------------------------

# The document currently carries a stray "_GoBack" bookmark further down in
# the file (in an otherwise empty paragraph). Word keeps only one "_GoBack"
# bookmark at a time (it marks the most recent edit position), so remove the
# old one before adding the new one - otherwise later bookmark lookups would
# resolve to whichever one happens to come first.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$keywordsRange = $d.Content
$foundKw = $keywordsRange.Find.Execute("Poema, versos, palabras, literatura")
if ($foundKw) {
    $kwStart = $keywordsRange.Start
    $para = $d.Range($kwStart, $kwStart).Paragraphs(1)
    $paraStart = $para.Range.Start

    $rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr>'

    $newRuns = ''
    $newRuns += '<w:r>' + $rPr + '<w:t>&#x201C;</w:t></w:r>'
    $newRuns += '<w:proofErr w:type="spellStart"/>'
    $newRuns += '<w:r>' + $rPr + '<w:t>Poema</w:t></w:r>'
    $newRuns += '<w:r>' + $rPr + '<w:t>,versos,</w:t></w:r>'
    $newRuns += '<w:r>' + $rPr + '<w:t>palabras</w:t></w:r>'
    $newRuns += '<w:r>' + $rPr + '<w:t>,</w:t></w:r>'
    $newRuns += '<w:r>' + $rPr + '<w:t>literatura</w:t></w:r>'
    $newRuns += '<w:proofErr w:type="spellEnd"/>'
    $newRuns += '<w:r>' + $rPr + '<w:t>&#x201D;</w:t></w:r>'
    $newRuns += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>'
    $pkg += '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">'
    $pkg += '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">'
    $pkg += '<pkg:xmlData>'
    $pkg += '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $pkg += '<w:p>' + $newRuns + '</w:p>'
    $pkg += '</w:body></w:document>'
    $pkg += '</pkg:xmlData></pkg:part></pkg:package>'

    # Inserting at the exact paragraph-start offset merges the fragment's
    # runs into the existing paragraph instead of splitting off a new one.
    $insertionPoint = $d.Range($paraStart, $paraStart)
    $insertionPoint.InsertXML($pkg)

    # The original three runs ("Poema" / ", versos, palabras" / ", literatura")
    # are still present right after what we just inserted; remove them.
    $para2 = $d.Range($paraStart, $paraStart).Paragraphs(1)
    $staleRange = $d.Range($paraStart, $para2.Range.End)
    $staleRange.Find.Execute("Poema, versos, palabras, literatura") | Out-Null
    $staleRange.Delete()
}
